# Corrige a ordem dos jogos (linhas trocadas) na base da Mexico Liga MX Femenil.
# Para cada par de linhas abaixo, os dados das colunas B:AD (id, times, odds, etc.)
# estavam atribuidos a linha errada; a coluna A (numero sequencial) permanece fixa
# por linha e o conteudo de B:AD e trocado entre as duas linhas do par.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(28, 29),
    @(71, 72),
    @(232, 233),
    @(245, 246),
    @(263, 265),
    @(271, 272),
    @(310, 311)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B" + $r1 + ":AD" + $r1)
    $range2 = $ws.Range("B" + $r2 + ":AD" + $r2)

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
